$wb = $excel.ActiveWorkbook

# --- Add the new "setup" worksheet at the end of the tab strip ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "setup"

# Column A source-code snippets, column B their Russian explanations
$codeLines = @(
    '@pytest.mark.add_to_basket',
    'class TestAddToBasketFromProductPage(object):',
    '    @pytest.fixture(scope="function", autouse=True)',
    '    def setup(self):',
    '        self.product = ProductFactory(title="Best book created by robot")',
    '        self.link = self.product.link',
    '        yield',
    '        self.product.delete()',
    '    def test_guest_cant_see_success_message(self, browser):',
    '        page = ProductPage(browser, self.link)'
)

$explainLines = @(
    'Метка добавить в карзину(вызвать можно через консоль данный метод)',
    'Название класса',
    'фикстура с автозапуском для тестов',
    'добавление сетапа',
    'Добавление продукта с названием таким то',
    'установка ссылки для продукта',
    'По завершению тест',
    'удаляем продукт из бд или функцией',
    'начало теста передаем селф и бруезер',
    'создание страницы Продукта'
)

# Column A width approximating the source workbook's 55.21875 chars
$ws.Columns.Item(1).ColumnWidth = 54.33

# Shared-string indices must come out in the same order as the source
# workbook: all of column A first (39-48), then all of column B (49-58).
for ($i = 0; $i -lt $codeLines.Length; $i++) {
    $row = $i + 1
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $codeLines[$i]
    $cellA.Font.Name = "Arial"
    $cellA.Font.Size = 10
    $cellA.Font.Color = 0
    $cellA.HorizontalAlignment = -4131
    $cellA.VerticalAlignment = -4108
    $cellA.WrapText = $true
}

for ($i = 0; $i -lt $explainLines.Length; $i++) {
    $row = $i + 1
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $explainLines[$i]
}

# Row 5 holds the longest snippet, which wraps onto two lines in the
# original workbook (row height 26.4pt instead of the 14.4pt default).
$ws.Rows.Item(5).RowHeight = 26.4

$ws.Range("B11").Select()

# Activate the new sheet so it becomes the selected/active tab, matching
# the tab-selection state recorded in the committed workbook.
$ws.Activate()
